{"js": "// Change every \"Justify\" aligned paragraph in the document body to\n// \"Align Left\" (mirrors the XML edit w:jc val=\"both\" -> w:jc val=\"left\"\n// that was applied to all the code-block paragraphs in the doc).\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/alignment\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  if (para.alignment === Word.Alignment.justified) {\n    para.alignment = Word.Alignment.left;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Change every \"Justify\" aligned paragraph in the document body to\n# \"Align Left\" (mirrors the XML edit w:jc val=\"both\" -> w:jc val=\"left\"\n# that was applied to all the code-block paragraphs in the doc).\n$d = $word.ActiveDocument\n\n$wdAlignParagraphLeft = 0\n$wdAlignParagraphJustify = 3\n\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.ParagraphFormat.Alignment -eq $wdAlignParagraphJustify) {\n        $p.Range.ParagraphFormat.Alignment = $wdAlignParagraphLeft\n    }\n}\n"}
